# Update the "Configuration" sheet values
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Configuration")

$ws1.Range("B1").Value = 300
$ws1.Range("B2").Value = 100
$ws1.Range("B6").Value = 50
$ws1.Range("B8").Value = 1

# New configuration rows (copy the formatting used by the other "label"
# cells in column A down onto the newly-added rows)
$ws1.Range("A11").Value = "SAVED_DETAILED_RESULTS"
$ws1.Range("B11").Value = 0
$ws1.Range("A10").Copy()
$ws1.Range("A11").PasteSpecial(-4122)

# Add the new "marketQuota" worksheet after the last sheet ("buyers")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws4.Name = "marketQuota"

$ws4.Range("A1").Value = "Aliexpress"
$ws4.Range("B1").Value = 92.270531400966178
$ws4.Range("A2").Value = "Banggood"
$ws4.Range("B2").Value = 10.144927536231885
$ws4.Range("A3").Value = "Wish"
$ws4.Range("B3").Value = 41.304347826086953
$ws4.Range("A4").Value = "Others"
$ws4.Range("B4").Value = 14.251207729468598
$ws4.Range("A5").Value = "Lightinthebox"
$ws4.Range("B5").Value = 11.594202898550725
$ws4.Range("A6").Value = "Alibaba"
$ws4.Range("B6").Value = 42.512077294685987

# NOTE: Excel's ColumnWidth (in "characters") is stored on a pixel grid, so
# it can only land on specific quantized values; the inputs below are chosen
# to land as close as possible to the target widths (14.5 and ~16.332 chars).
$ws4.Columns.Item(1).ColumnWidth = 13.666666666666666
$ws4.Columns.Item(2).ColumnWidth = 15.5
$ws4.Range("B2").Select()

# Finish configuration-row additions (kept after sheet4 writes so that
# shared-string insertion order matches: SAVED_DETAILED_RESULTS, Banggood,
# Lightinthebox, MARKET_QUOTA)
$ws1.Range("A12").Value = "MARKET_QUOTA"
$ws1.Range("B12").Value = 1
$ws1.Range("A10").Copy()
$ws1.Range("A12").PasteSpecial(-4122)

# Same quantization caveat as above -- closest achievable value to the
# target width of ~23.664 characters.
$ws1.Columns.Item(1).ColumnWidth = 22.833333333333332
$ws1.Range("B7").Select()
